$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'249.87"
$ws.Range("E2").Value = "'1.77%"
$ws.Range("D3").Value = "'28.44"
$ws.Range("E3").Value = "'-3.32%"
$ws.Range("D4").Value = "'5.262"
$ws.Range("E4").Value = "'1.88%"
$ws.Range("D5").Value = "'0.05753"
$ws.Range("E5").Value = "'-0.05%"
$ws.Range("D6").Value = "'6.661"
$ws.Range("E6").Value = "'1.24%"
$ws.Range("D7").Value = "'3.233"
$ws.Range("E7").Value = "'3.45%"
$ws.Range("D8").Value = "'0.8617"
$ws.Range("E8").Value = "'0.36%"
$ws.Range("D9").Value = "'0.9189"
$ws.Range("E9").Value = "'7.18%"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "'0.01043"
$ws.Range("E10").Value = "'1,637.67%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1404"
$ws.Range("E11").Value = "'2.88%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.07166"
$ws.Range("E12").Value = "'1.82%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03134"
$ws.Range("E13").Value = "'2.86%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09250"
$ws.Range("E14").Value = "'-1.26%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001526"
$ws.Range("E15").Value = "'-0.35%"
$ws.Range("D16").Value = "'0.005956"
$ws.Range("E16").Value = "'-0.14%"
$ws.Range("D17").Value = "'3.501"
$ws.Range("E17").Value = "'0.26%"
$ws.Range("D18").Value = "'2.255"
$ws.Range("E18").Value = "'2.93%"
$ws.Range("D19").Value = "'0.3127"
$ws.Range("E19").Value = "'-2.30%"
$ws.Range("D20").Value = "'0.03379"
$ws.Range("E20").Value = "'1.84%"
$ws.Range("D21").Value = "'0.1313"
$ws.Range("E21").Value = "'2.42%"
$ws.Range("D22").Value = "'3.539"
$ws.Range("E22").Value = "'0.90%"
$ws.Range("D23").Value = "'0.04173"
$ws.Range("E23").Value = "'0.58%"
$ws.Range("D24").Value = "'0.1378"
$ws.Range("E24").Value = "'-1.54%"
$ws.Range("D25").Value = "'0.005042"
$ws.Range("E25").Value = "'22.05%"
$ws.Range("D26").Value = "'0.001218"
$ws.Range("E26").Value = "'-0.72%"
$ws.Range("E27").Value = "'-0.90%"
$ws.Range("D28").Value = "'0.0001937"
$ws.Range("E28").Value = "'33.67%"
$ws.Range("D41").Value = "'0.005676"
$ws.Range("E41").Value = "'-1.30%"
$ws.Range("D42").Value = "'0.1082"
$ws.Range("E42").Value = "'1.12%"
$ws.Range("D43").Value = "'0.002199"
$ws.Range("E43").Value = "'-0.07%"
$ws.Range("D44").Value = "'0.009757"
$ws.Range("E44").Value = "'16.25%"
$ws.Range("D45").Value = "'0.00005289"
$ws.Range("E45").Value = "'0.29%"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("E46").Value = "'-0.07%"
$ws.Range("D47").Value = "'0.08495"
$ws.Range("E47").Value = "'46.45%"
$ws.Range("E48").Value = "'-11.10%"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("E49").Value = "'-0.07%"
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("E50").Value = "'-0.07%"
